# Apply the cryptos list update (prices, 1h volume %, and two name/link/price swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D hold price text that looks numeric (e.g. "67.019.50" / "580.39").
# Force text format first so Excel does not silently coerce these into floating-point
# numbers (which would lose the literal string representation / thousands-dot format).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.019.50"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.118.49"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "580.39"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "173.45"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("B9").Value = "Toncoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D9").Value = "6.40"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "0.155"
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("D11").Value = "0.480"
$ws.Range("E11").Value = "  -0.83%  "
$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "37.15"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "3.641.33"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").Value = "67.024.83"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "7.12"
$ws.Range("E17").Value = "  -1.03%  "
$ws.Range("D18").Value = "3.119.33"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "16.34"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "491.06"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("E21").Value = "  +5.28%  "
$ws.Range("D22").Value = "0.706"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").Value = "84.15"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("E25").Value = "  -4.24%  "
$ws.Range("D26").Value = "10.42"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("E29").Value = "  -1.83%  "
$ws.Range("D30").Value = "2.68"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "28.58"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "0.0₃0945"
$ws.Range("E33").Value = "  -7.18%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "5.88"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").Value = "0.974"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "47.27"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("D39").Value = "0.310"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "8.52"
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "2.829.31"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "383.49"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -7.71%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "135.45"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D48").Value = "24.99"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "2.21"
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("E50").Value = "  -0.98%  "
$ws.Range("D51").Value = "6.75"
$ws.Range("E51").Value = "  -0.93%  "

# Reset those cells back to the default (unstyled) look now that the text value is set,
# so no lingering quote-prefix / text-number-format styling is left on the cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
